$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.855689333333333
$ws.Range("I2").Value = 0.1513312545414155
$ws.Range("J2").Value = 0.1513312545414155
$ws.Range("O2").Value = 0.1284798252444435
$ws.Range("P2").Value = 0.1284798252444435
$ws.Range("Q2").Value = 0.09479212225999999
$ws.Range("R2").Value = 0.8531291003399999
$ws.Range("S2").Value = 0.01944301313750346
$ws.Range("T2").Value = 0.01944301313750346
$ws.Range("G3").Value = 3.855689333333333
$ws.Range("I3").Value = 0.1513312545414155
$ws.Range("J3").Value = 0.1513312545414155
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.166768
$ws.Range("N3").Value = 0.500304
$ws.Range("O3").Value = 0.8715201747555565
$ws.Range("P3").Value = 0.8715201747555564
$ws.Range("Q3").Value = 0.6430055987413332
$ws.Range("R3").Value = 5.787050388671999
$ws.Range("S3").Value = 0.131888241403912
$ws.Range("T3").Value = 0.131888241403912
$ws.Range("G4").Value = 5.360192666666666
$ws.Range("H4").Value = 16.080578
$ws.Range("I4").Value = 0.2103812342497758
$ws.Range("J4").Value = 0.2103812342497758
$ws.Range("O4").Value = 0.1284798252444435
$ws.Range("P4").Value = 0.1284798252444435
$ws.Range("Q4").Value = 0.13178033671
$ws.Range("R4").Value = 1.18602303039
$ws.Range("S4").Value = 0.02702974421112153
$ws.Range("T4").Value = 0.02702974421112153
$ws.Range("G5").Value = 5.360192666666666
$ws.Range("H5").Value = 16.080578
$ws.Range("I5").Value = 0.2103812342497758
$ws.Range("J5").Value = 0.2103812342497758
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.166768
$ws.Range("N5").Value = 0.500304
$ws.Range("O5").Value = 0.8715201747555565
$ws.Range("P5").Value = 0.8715201747555564
$ws.Range("Q5").Value = 0.8939086106346666
$ws.Range("R5").Value = 8.045177495711998
$ws.Range("S5").Value = 0.1833514900386543
$ws.Range("T5").Value = 0.1833514900386543
$ws.Range("G6").Value = 16.26259133333333
$ws.Range("H6").Value = 48.787774
$ws.Range("I6").Value = 0.6382875112088087
$ws.Range("J6").Value = 0.6382875112088088
$ws.Range("O6").Value = 0.1284798252444435
$ws.Range("P6").Value = 0.1284798252444435
$ws.Range("Q6").Value = 0.39981580793
$ws.Range("R6").Value = 3.59834227137
$ws.Range("S6").Value = 0.08200706789581852
$ws.Range("T6").Value = 0.08200706789581853
$ws.Range("G7").Value = 16.26259133333333
$ws.Range("H7").Value = 48.787774
$ws.Range("I7").Value = 0.6382875112088087
$ws.Range("J7").Value = 0.6382875112088088
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.166768
$ws.Range("N7").Value = 0.500304
$ws.Range("O7").Value = 0.8715201747555565
$ws.Range("P7").Value = 0.8715201747555564
$ws.Range("Q7").Value = 2.712079831477333
$ws.Range("R7").Value = 24.408718483296
$ws.Range("S7").Value = 0.5562804433129902
$ws.Range("T7").Value = 0.5562804433129902
